$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (Price) and hourly change (Volume(1h)) columns
# with refreshed data per the GitHub Actions scheduled run.

$ws.Range("D2").Value = "27.015.13"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.874.20"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'305.77"
$ws.Range("E5").Value = "  -0.01%  "
$ws.Range("D6").Value = "'0.9990"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").Value = "'0.5062"
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("D8").Value = "'0.3660"
$ws.Range("E8").Value = "  -2.03%  "
$ws.Range("D9").Value = "'0.07196"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'0.8953"
$ws.Range("E10").Value = "  +1.42%  "
$ws.Range("D11").Value = "'20.74"
$ws.Range("E11").Value = "  +0.87%  "
$ws.Range("D12").Value = "1.864.41"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "'0.07520"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "'95.16"
$ws.Range("E14").Value = "  +6.91%  "
$ws.Range("D15").Value = "'5.240"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'0.000008538"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "'0.9997"
$ws.Range("D20").Value = "27.053.22"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'5.027"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").Value = "2.109.84"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Value = "'10.42"
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'6.420"
$ws.Range("E24").Value = "  -0.47%  "
$ws.Range("D25").Value = "'148.28"
$ws.Range("E25").Value = "  +0.87%  "
$ws.Range("D26").Value = "'1.778"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("D27").Value = "'17.91"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'2.078"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'113.36"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("D30").Value = "'4.705"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("D32").Value = "'0.09167"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("D33").Value = "'0.05141"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").Value = "'0.7549"
$ws.Range("E34").Value = "  +4.52%  "
$ws.Range("D35").Value = "'2.993"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").Value = "'1.161"
$ws.Range("E36").Value = "  +1.08%  "
$ws.Range("D37").Value = "'3.229"
$ws.Range("E37").Value = "  +6.32%  "
$ws.Range("D38").Value = "'2.594"
$ws.Range("E38").Value = "  +5.85%  "
$ws.Range("D39").Value = "'0.5640"
$ws.Range("E39").Value = "  +6.96%  "
$ws.Range("D40").Value = "'0.02002"
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").Value = "'1.071"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "'6.604"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("D43").Value = "'116.12"
$ws.Range("E43").Value = "  +1.15%  "
$ws.Range("D44").Value = "'8.558"
$ws.Range("E44").Value = "  +3.98%  "
$ws.Range("D45").Value = "'0.1477"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").Value = "'0.4731"
$ws.Range("E46").Value = "  +2.86%  "
$ws.Range("D47").Value = "'0.9986"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "'10.14"
$ws.Range("E48").Value = "  +2.00%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "'36.90"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "'63.27"
$ws.Range("E51").Value = "  -0.83%  "
